# Add a new header row ("date" / "time") above the existing data table,
# shifting all existing rows down by one, and make the new header bold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 1; everything that used to be in row 1 (and below)
# moves down to row 2 (and below).
$ws.Rows("1:1").Insert()

# Populate the new header row.
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "time"

# Make the new header bold (this introduces a new bold font + cell style).
$ws.Range("A1:B1").Font.Bold = $true

# Select the whole header row, matching the state left behind after the edit.
$ws.Rows("1:1").Select()
